$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new sprint log entry (SF-15 / preGameScreen takes the name ...)
$ws.Range("A5").Value = 42830
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = "SF-15"
$ws.Range("D5").Value = "preGameScreen takes the name and places it into the GameBoard as playerName"

# Row 6: new sprint log entry (SF-2 / preGameScreen gives the option ...)
$ws.Range("A6").Value = 42830
$ws.Range("B6").Value = 0.5
$ws.Range("C6").Value = "SF-2"
$ws.Range("D6").Value = "preGameScreen gives the option of a timed game or untimed game"

# Match the font size used elsewhere in the log table for the new rows
$ws.Range("C5:D6").Font.Size = 12

# Selection moved to A7 (next empty row) after entering the new data
$ws.Range("A7").Select()
